$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove the now-obsolete data rows (old rows 4,5,6) ---
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(4).Delete()

# --- 2. Preserve the date number format (style s="1", numFmtId 14) that
#        currently lives on column C (old "Date" column) by copying its
#        format to the new "Date" column (B) before we overwrite column C. ---
$ws.Range("C2").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Header row ---
$ws.Range("A1").Value = "Bill No"
$ws.Range("B1").Value = "Date"
$ws.Range("C1").Value = "Cashier"
$ws.Range("D1").Value = "KOT"
$ws.Range("E1").Value = "Price (₹)"
$ws.Range("F1").Value = "SGST (₹)"
$ws.Range("G1").Value = "CGST (₹)"
$ws.Range("H1").Value = "Tax (₹)"
$ws.Range("I1").Value = "Food Items"

# --- 4. Data row 2 ---
$ws.Range("A2").Value = 432
$ws.Range("B2").Value = 45903.00011574074
$ws.Range("C2").ClearFormats()
$ws.Range("C2").Value = "Ajay Francis Anchan"
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 99
$ws.Range("F2").Value = 1.93
$ws.Range("G2").Value = 1.93
$ws.Range("H2").Value = 3.85
$ws.Range("I2").Value = "Chicken Burger (x1), Strawberry Lassi (x1), Veg Wrap (x1)"

# --- 5. Data row 3 ---
$ws.Range("A3").Value = 433
$ws.Range("B3").Value = 45903.00011574074
$ws.Range("C3").ClearFormats()
$ws.Range("C3").Value = "Ajay Francis Anchan"
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 500
$ws.Range("F3").Value = 7.1
$ws.Range("G3").Value = 7.1
$ws.Range("H3").Value = 16
$ws.Range("I3").Value = "Chicken Burger (x1), Chicken Cheese Burger (x1), Mango Lassi (x1), Strawberry Lassi (x1), Chicken Cheese Pops (x1), Veg Cheese Pops (x3)"
